$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table ---
# Update Intel 23.90.0.2 row (row 3): Critical Minutes / Good Roaming %
$ws.Cells.Item(3, 3).Value = 662
$ws.Cells.Item(3, 4).Value = 97.5

# Remove the MediaTek MT7921 ...3.0.1.1297 row (row 4); rows below shift up
$ws.Rows(4).Delete()

# Totals row is now row 4 (was row 5): Client Count / Critical Minutes totals
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 662

# --- Good Drivers table ---
# Insert a new data row right after the header row (now row 11), for the
# new "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3" entry
$ws.Rows(12).Insert()
$ws.Cells.Item(12, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Cells.Item(12, 2).Value = 11128
$ws.Cells.Item(12, 2).NumberFormat = $ws.Cells.Item(13, 2).NumberFormat
$ws.Cells.Item(12, 4).Value = 100

# Row 13 (was row 12): Intel 23.100.0.4 total samples updated
$ws.Cells.Item(13, 2).Value = 486214

# Row 14: MediaTek ...3.0.1.1255 -> Intel 22.80.0.9
$ws.Cells.Item(14, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Cells.Item(14, 2).Value = 79953
$ws.Cells.Item(14, 5).Value = "'2021-08-18"

# Row 15: MediaTek ...3.0.1.1216 -> Intel 22.50.1.1
$ws.Cells.Item(15, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Cells.Item(15, 2).Value = 35355
$ws.Cells.Item(15, 5).Value = "'2021-04-27"

# Row 16: Intel 22.80.0.9 -> Intel 21.110.3.2
$ws.Cells.Item(16, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Cells.Item(16, 2).Value = 65425
$ws.Cells.Item(16, 4).Value = 100
$ws.Cells.Item(16, 5).Value = "'2020-08-05"

# Row 17: Intel 22.50.1.1 -> Intel 21.70.0.6
$ws.Cells.Item(17, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Cells.Item(17, 2).Value = 117653
$ws.Cells.Item(17, 5).Value = "'2020-01-06"

# Row 18: Intel 21.110.3.2 -> Intel 21.60.2.1
$ws.Cells.Item(18, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Cells.Item(18, 2).Value = 56018
$ws.Cells.Item(18, 5).Value = "'2019-12-14"

# Old rows 19 and 20 (Intel 21.70.0.6 / Intel 21.60.2.1) are now duplicated
# by the renamed rows above; remove them so the table ends at row 18.
$ws.Rows(19).Delete()
$ws.Rows(19).Delete()

# --- Column width ---
# The stored OOXML <col width> differs from the COM ColumnWidth property by
# the standard Calibri-11 5/6-character padding; 43 + 1/6 round-trips to
# exactly width="44" in the saved file.
$ws.Columns(1).ColumnWidth = 43.16666667
